# Update "想去人数" (number of people interested) counts for several events.
# The same rows exist both in the "展览" sheet and the combined "全部类型"
# sheet, so the same values need to be applied in both places.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    "F2"  = 1188
    "F6"  = 181
    "F10" = 5537
    "F11" = 4921
    "F15" = 53
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
